{"js": "// Update the date line and the 25 \"two-digit number divided by one-digit\n// number\" problems in the practice table. Values are replaced in place\n// (table row/column position) so runs keep their original formatting\n// (font, size, paragraph alignment) and duplicate problem strings (e.g.\n// \"78\u00f77=\" appearing twice) are each mapped to their own distinct result.\n\n// 1. Update the date heading paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2025-08-05 Tuesday\", \"Replace\");\n\n// 2. Update the practice-problem table, cell by cell, by (row, column)\n//    position so duplicate values resolve to the correct replacement.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indices in the 20-row table that hold visible problems (the other\n// rows are blank spacer rows). Each entry has 5 new values, one per\n// column, matching the document's left-to-right order.\nconst newValuesByRow = {\n  0: [\"94\u00f72=\", \"37\u00f75=\", \"77\u00f79=\", \"19\u00f73=\", \"31\u00f76=\"],\n  4: [\"29\u00f79=\", \"17\u00f73=\", \"73\u00f73=\", \"10\u00f76=\", \"43\u00f76=\"],\n  8: [\"34\u00f74=\", \"63\u00f73=\", \"59\u00f73=\", \"60\u00f75=\", \"80\u00f79=\"],\n  12: [\"83\u00f79=\", \"49\u00f75=\", \"40\u00f79=\", \"63\u00f76=\", \"90\u00f78=\"],\n  16: [\"25\u00f78=\", \"33\u00f73=\", \"74\u00f77=\", \"12\u00f72=\", \"64\u00f77=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = newValuesByRow[rowIndex];\n  for (let colIndex = 0; colIndex < newValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    const cellParagraph = cell.body.paragraphs.getFirst();\n    cellParagraph.getRange().insertText(newValues[colIndex], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"two-digit number divided by one-digit\n# number\" problems in the practice table. Values are assigned by table\n# (row, column) position so each run keeps its original formatting (font,\n# size, paragraph alignment) and duplicate problem strings (e.g. \"78\u00f77=\"\n# appearing twice) each resolve to their own distinct replacement.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date heading paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2025-08-05 Tuesday\"\n\n# 2. Update the practice-problem table, cell by cell, by (row, column)\n#    position so duplicate values resolve to the correct replacement.\n$t = $d.Tables.Item(1)\n\n# 1-based row indices in the 20-row table that hold visible problems (the\n# other rows are blank spacer rows). Each entry has 5 new values, one per\n# column, matching the document's left-to-right order.\n$newValuesByRow = @{\n  1  = @(\"94\u00f72=\", \"37\u00f75=\", \"77\u00f79=\", \"19\u00f73=\", \"31\u00f76=\")\n  5  = @(\"29\u00f79=\", \"17\u00f73=\", \"73\u00f73=\", \"10\u00f76=\", \"43\u00f76=\")\n  9  = @(\"34\u00f74=\", \"63\u00f73=\", \"59\u00f73=\", \"60\u00f75=\", \"80\u00f79=\")\n  13 = @(\"83\u00f79=\", \"49\u00f75=\", \"40\u00f79=\", \"63\u00f76=\", \"90\u00f78=\")\n  17 = @(\"25\u00f78=\", \"33\u00f73=\", \"74\u00f77=\", \"12\u00f72=\", \"64\u00f77=\")\n}\n\nforeach ($rowIndex in $newValuesByRow.Keys) {\n  $values = $newValuesByRow[$rowIndex]\n  for ($col = 1; $col -le $values.Length; $col++) {\n    $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n  }\n}\n"}
